$p = $ppt.ActivePresentation

# --- Slide 4: "Content Placeholder 3" - update the two Edge bullet paragraphs ---
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(3)
$tr4 = $shape4.TextFrame.TextRange

# Paragraph 5: "Graph edge (to a graph node);" -> "Graph edge  -> ;"
# (go through a disjoint placeholder first so the host doesn't keep the
#  shared prefix/suffix as a separate run - we want a single clean run)
$para5 = $tr4.Paragraphs(5, 1)
$para5.Text = "zzzzPLACEHOLDERzzzz"
$para5b = $tr4.Paragraphs(5, 1)
$para5b.Text = "Graph edge  -> ;"

# Paragraph 6: "Virtual edge (to a Time-Series node)" -> "Virtual edge  -> "
$para6 = $tr4.Paragraphs(6, 1)
$para6.Text = "zzzzPLACEHOLDERzzzz"
$para6b = $tr4.Paragraphs(6, 1)
$para6b.Text = "Virtual edge  -> "

# --- Slide 5: "Content Placeholder 2" - add a new bullet after the last one ---
$s5 = $p.Slides.Item(5)
$shape5 = $s5.Shapes.Item(2)
$tr5 = $shape5.TextFrame.TextRange
$lastCount = $tr5.Paragraphs().Count
$lastPara = $tr5.Paragraphs($lastCount, 1)
$lastPara.InsertAfter("`rOutgoing edges and properties are stored in a “fat” representation”.") | Out-Null
